# Refresh the rolling 10-day forecast window: the series advances by three days
# (03.02.2026 -> 06.02.2026 ... 10.02.2026 -> 13.02.2026) and the updated NRG /
# PCSun / Ulmeni model predictions are written into the "Prediction" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Prediction" (column C) values for data rows 2..170, in row order.
$newPredictions = @(
    0, 0, 0.178, 0.042, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0,
    0, 0, 0, 0, 0.058, 0.182, 0.251, 0.357, 0.502, 0.471, 0.486, 0.307, 0.084, 0, 0,
    0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0.098, 0.194,
    0.275, 0.328, 0.302, 0.281, 0.2, 0.119, 0.061, 0, 0, 0, 0, 0, 0, 0, 0,
    0, 0, 0, 0, 0, 0, 0, 0.094, 0.266, 0.472, 0.482, 0.482, 0.505, 0.348, 0.201,
    0.064, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0,
    0, 0.242, 0.661, 0.974, 1.252, 1.243, 1.038, 0.582, 0.346, 0.089, 0, 0, 0, 0, 0,
    0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0.178, 0.478, 0.727, 1.121, 1.33,
    1.219, 1.101, 0.55, 0.193, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0,
    0, 0, 0, 0.015, 0.193, 0.517, 0.843, 1.227, 1.274, 1.112, 0.806, 0.484, 0.183, 0, 0,
    0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0.093, 0.238,
    0.329, 0.422, 0.446, 0.336
)

$firstRow = 2
$lastRow = 170

# The window now starts on 06.02.2026 (Excel serial 46059), keeping the same
# hour-of-day (column B) pattern as before, just shifted three days later.
$currentDate = Get-Date -Year 2026 -Month 2 -Day 6 -Hour 0 -Minute 0 -Second 0
$previousHour = 0

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $hour = $ws.Cells.Item($row, 2).Value()
    if ($hour -le $previousHour) {
        $currentDate = $currentDate.AddDays(1)
    }
    $previousHour = $hour

    $ws.Cells.Item($row, 1).Value = $currentDate
    $ws.Cells.Item($row, 3).Value = $newPredictions[$row - $firstRow]
    $ws.Cells.Item($row, 4).Value = $currentDate.ToString("dd.MM.yyyy") + $hour
}
